$wb = $excel.ActiveWorkbook

# Sheets "展览" and "全部类型" both list the same three events in rows 2-4
# and need their "想去人数" (column F) counts refreshed.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 3205
    $ws.Range("F3").Value = 50
    $ws.Range("F4").Value = 1036
}
